$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "punktindeks_maned" (2nd sheet) - monthly point-index figures.
# October ("okt", column N) values added for several measurement points,
# and a handful of previously-present provisional values removed.
# ---------------------------------------------------------------------------
$wsMonth = $wb.Worksheets.Item(2)

# Row 4 (Breivika, E8): remove provisional feb/mar/apr values.
$wsMonth.Range("F4:H4").ClearContents()

# New October ("okt", column N) figures.
$wsMonth.Range("N13").Value = 1.78
$wsMonth.Range("N16").Value = 3.49
$wsMonth.Range("N19").Value = 5.67
$wsMonth.Range("N25").Value = -2.89
$wsMonth.Range("N34").Value = 5.88
$wsMonth.Range("N40").Value = 2.24
$wsMonth.Range("N49").Value = 2.43
$wsMonth.Range("N58").Value = -4.11
$wsMonth.Range("N61").Value = 7.25
$wsMonth.Range("N64").Value = 2.87
$wsMonth.Range("N73").Value = -0.33
$wsMonth.Range("N76").Value = 9.93
$wsMonth.Range("N79").Value = 8.1

# Row 55 (Prestvannet, kv. 21200): remove provisional jan/feb/mai/jul values.
$wsMonth.Range("E55").ClearContents()
$wsMonth.Range("F55").ClearContents()
$wsMonth.Range("I55").ClearContents()
$wsMonth.Range("K55").ClearContents()

# Row 82 (Tverrforbindelsen, E8): remove provisional jul value.
$wsMonth.Range("K82").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "byindeks_aarlig" (3rd sheet) - annual city index, recomputed
# standard errors plus the rolling 2024-2025 ("jan-okt") window extended
# from September to October.
# ---------------------------------------------------------------------------
$wsYearly = $wb.Worksheets.Item(3)

$wsYearly.Range("J3").Value = 1.18240265755164
$wsYearly.Range("J4").Value = 2.077955497173004

$wsYearly.Range("C5").Value = 10
$wsYearly.Range("F5").Value = "jan-okt"
$wsYearly.Range("G5").Value = 18
$wsYearly.Range("H5").Value = 1.014
$wsYearly.Range("I5").Value = 1.4
$wsYearly.Range("J5").Value = 1.116766126251205
$wsYearly.Range("L5").Value = 3.6

$wsYearly.Range("J6").Value = 2.306349631338904
$wsYearly.Range("J7").Value = 1.397320590274238
$wsYearly.Range("J8").Value = 2.368639008242137

# ---------------------------------------------------------------------------
# Sheet "by_glid_indeks" (4th sheet) - rolling window index. A new 12-month
# window (nov 2024 - okt 2025) is inserted as the new row 24 (pushing the
# existing rolling windows down by one row), and a new 24-month window
# (nov 2023 - okt 2025) is appended as the new last row (35).
# ---------------------------------------------------------------------------
$wsRoll = $wb.Worksheets.Item(4)

$wsRoll.Rows(24).Insert()

$wsRoll.Cells.Item(24, 1).Value = 0.9568554091875735
$wsRoll.Cells.Item(24, 2).Value = -4.314459081242649
$wsRoll.Cells.Item(24, 3).Value = 14
$wsRoll.Cells.Item(24, 4).Value = 9.461194469994499
$wsRoll.Cells.Item(24, 5).Value = 10.62561012946104
$wsRoll.Cells.Item(24, 6).Value = 3.422293990449863
$wsRoll.Cells.Item(24, 7).Value = -11
$wsRoll.Cells.Item(24, 8).Value = 2.4
$wsRoll.Cells.Item(24, 9).Value = "2019 - (nov 2024 - okt 2025)"
$wsRoll.Cells.Item(24, 10).Value = 45931
$wsRoll.Cells.Item(24, 11).Value = 10
$wsRoll.Cells.Item(24, 12).Value = 2025
$wsRoll.Cells.Item(24, 13).Value = "12_months"

$wsRoll.Cells.Item(35, 1).Value = 0.9504647415344348
$wsRoll.Cells.Item(35, 2).Value = -4.953525846556516
$wsRoll.Cells.Item(35, 3).Value = 14
$wsRoll.Cells.Item(35, 4).Value = 9.461194469994499
$wsRoll.Cells.Item(35, 5).Value = 9.445356864131334
$wsRoll.Cells.Item(35, 6).Value = 3.066579836549348
$wsRoll.Cells.Item(35, 7).Value = -11
$wsRoll.Cells.Item(35, 8).Value = 1.1
$wsRoll.Cells.Item(35, 9).Value = "2019 - (nov 2023 - okt 2025)"
$wsRoll.Cells.Item(35, 10).Value = 45931
$wsRoll.Cells.Item(35, 11).Value = 10
$wsRoll.Cells.Item(35, 12).Value = 2025
$wsRoll.Cells.Item(35, 13).Value = "24_months"
